$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.409.48'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.869.62'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.98'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7050'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07911'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3140'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.55'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07888'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.874.58'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.70%  '
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.192'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.14%  '
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.82'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008394'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.395.21'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '254.38'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.128.72'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.645'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.81%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1559'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.009'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.23'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.80'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.508'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.332'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.263'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.215'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05298'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.893'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7517'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.176'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.713'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01885'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.284.69'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.764'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8933'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.18'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.024'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -7.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.17'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.12%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.030.40'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.799'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.588'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.20%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4312'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.09%  '
